# Natmi following Dr Hou advice
# Updates recomputed ligand/receptor-expressing cell counts and derived
# specificity/weight statistics for the Serpine1-Lrp1 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.145195333333334
$ws.Range("H2").Value = 12.435586
$ws.Range("I2").Value = 0.04167134630420959
$ws.Range("J2").Value = 0.04167134630420959
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 97.97499818747755
$ws.Range("R2").Value = 881.774983687298
$ws.Range("S2").Value = 0.002845254763220851
$ws.Range("T2").Value = 0.002845254763220851

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.145195333333334
$ws.Range("H3").Value = 12.435586
$ws.Range("I3").Value = 0.04167134630420959
$ws.Range("J3").Value = 0.04167134630420959
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 751.4757937424406
$ws.Range("R3").Value = 6763.282143681965
$ws.Range("S3").Value = 0.02182332351259109
$ws.Range("T3").Value = 0.02182332351259109

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.145195333333334
$ws.Range("H4").Value = 12.435586
$ws.Range("I4").Value = 0.04167134630420959
$ws.Range("J4").Value = 0.04167134630420959
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 460.5334661917171
$ws.Range("R4").Value = 4144.801195725454
$ws.Range("S4").Value = 0.01337417772437447
$ws.Range("T4").Value = 0.01337417772437448

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.145195333333334
$ws.Range("H5").Value = 12.435586
$ws.Range("I5").Value = 0.04167134630420959
$ws.Range("J5").Value = 0.04167134630420959
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 124.9487859770164
$ws.Range("R5").Value = 1124.539073793148
$ws.Range("S5").Value = 0.003628590304023169
$ws.Range("T5").Value = 0.00362859030402317

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.36098766666667
$ws.Range("H6").Value = 190.082963
$ws.Range("I6").Value = 0.6369633869850008
$ws.Range("J6").Value = 0.6369633869850008
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 1497.587484449495
$ws.Range("R6").Value = 13478.28736004546
$ws.Range("S6").Value = 0.0434908701433839
$ws.Range("T6").Value = 0.0434908701433839

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.36098766666667
$ws.Range("H7").Value = 190.082963
$ws.Range("I7").Value = 0.6369633869850008
$ws.Range("J7").Value = 0.6369633869850008
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 11486.61152738117
$ws.Range("R7").Value = 103379.5037464306
$ws.Range("S7").Value = 0.333578328820281
$ws.Range("T7").Value = 0.3335783288202809

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 63.36098766666667
$ws.Range("H8").Value = 190.082963
$ws.Range("I8").Value = 0.6369633869850008
$ws.Range("J8").Value = 0.6369633869850008
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 7039.440345986261
$ws.Range("R8").Value = 63354.96311387635
$ws.Range("S8").Value = 0.2044297172274549
$ws.Range("T8").Value = 0.2044297172274549

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 63.36098766666667
$ws.Range("H9").Value = 190.082963
$ws.Range("I9").Value = 0.6369633869850008
$ws.Range("J9").Value = 0.6369633869850008
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 1909.892743435181
$ws.Range("R9").Value = 17189.03469091663
$ws.Range("S9").Value = 0.05546447079388096
$ws.Range("T9").Value = 0.05546447079388096

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.023468333333334
$ws.Range("H10").Value = 18.070405
$ws.Range("I10").Value = 0.06055348775782022
$ws.Range("J10").Value = 0.06055348775782022
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 142.3694787782406
$ws.Range("R10").Value = 1281.325309004165
$ws.Range("S10").Value = 0.004134498036488178
$ws.Range("T10").Value = 0.004134498036488178

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.023468333333334
$ws.Range("H11").Value = 18.070405
$ws.Range("I11").Value = 0.06055348775782022
$ws.Range("J11").Value = 0.06055348775782022
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 1091.984884397275
$ws.Range("R11").Value = 9827.86395957547
$ws.Range("S11").Value = 0.03171191886884492
$ws.Range("T11").Value = 0.03171191886884491

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.023468333333334
$ws.Range("H12").Value = 18.070405
$ws.Range("I12").Value = 0.06055348775782022
$ws.Range("J12").Value = 0.06055348775782022
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 669.2106226548661
$ws.Range("R12").Value = 6022.895603893795
$ws.Range("S12").Value = 0.01943429187988609
$ws.Range("T12").Value = 0.01943429187988609

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.023468333333334
$ws.Range("H13").Value = 18.070405
$ws.Range("I13").Value = 0.06055348775782022
$ws.Range("J13").Value = 0.06055348775782022
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 181.5656428947544
$ws.Range("R13").Value = 1634.09078605279
$ws.Range("S13").Value = 0.005272778972601034
$ws.Range("T13").Value = 0.005272778972601034

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.94386466666667
$ws.Range("H14").Value = 77.831594
$ws.Range("I14").Value = 0.2608117789529694
$ws.Range("J14").Value = 0.2608117789529694
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 613.2039359527157
$ws.Range("R14").Value = 5518.835423574442
$ws.Range("S14").Value = 0.01780782293311883
$ws.Range("T14").Value = 0.01780782293311883

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.94386466666667
$ws.Range("H15").Value = 77.831594
$ws.Range("I15").Value = 0.2608117789529694
$ws.Range("J15").Value = 0.2608117789529694
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 4703.321490389707
$ws.Range("R15").Value = 42329.89341350736
$ws.Range("S15").Value = 0.1365873755657871
$ws.Range("T15").Value = 0.1365873755657871

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.94386466666667
$ws.Range("H16").Value = 77.831594
$ws.Range("I16").Value = 0.2608117789529694
$ws.Range("J16").Value = 0.2608117789529694
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 2882.377538464729
$ws.Range("R16").Value = 25941.39784618256
$ws.Range("S16").Value = 0.0837060328903968
$ws.Range("T16").Value = 0.0837060328903968

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.94386466666667
$ws.Range("H17").Value = 77.831594
$ws.Range("I17").Value = 0.2608117789529694
$ws.Range("J17").Value = 0.2608117789529694
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 782.0269331060101
$ws.Range("R17").Value = 7038.242397954091
$ws.Range("S17").Value = 0.02271054756366672
$ws.Range("T17").Value = 0.02271054756366672
